# Applies the cryptocurrency price/volume refresh described by the commit:
# "Updated cryptos list on Thu Jan 11 03:53:37 UTC 2024 with GitHub Actions"
#
# Rows 31/32 also swap which coin (Filecoin / LidoDAOToken) occupies which
# row, which is handled below as a set of independent per-cell writes
# (same net effect as swapping the two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks like a plain number (e.g. "8.20", "10.50")
# as literal text, matching the workbook's existing convention of storing the
# Price column as text (so trailing zeros / exact formatting are preserved
# instead of being normalized away by Excel's automatic number detection).
function Set-TextCell($address, $value) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "46.435.19"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.602.47"
$ws.Range("E3").Value = "  +10.23%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextCell "D5" "305.74"
$ws.Range("E5").Value = "  +1.41%  "
Set-TextCell "D6" "100.13"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +5.49%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +13.40%  "
$ws.Range("E10").Value = "  +12.06%  "
Set-TextCell "D11" "0.0839"
$ws.Range("E11").Value = "  +4.99%  "
Set-TextCell "D12" "8.20"
$ws.Range("E12").Value = "  +14.87%  "
$ws.Range("D13").Value = "2.996.02"
$ws.Range("E13").Value = "  +10.12%  "
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "2.597.89"
$ws.Range("E15").Value = "  +9.89%  "
Set-TextCell "D16" "0.904"
$ws.Range("E16").Value = "  +10.95%  "
Set-TextCell "D17" "14.89"
$ws.Range("E17").Value = "  +9.35%  "
$ws.Range("D18").Value = "46.619.25"
$ws.Range("E18").Value = "  +1.30%  "
Set-TextCell "D19" "13.45"
$ws.Range("E19").Value = "  +5.28%  "
$ws.Range("E20").Value = "  +4.60%  "
$ws.Range("E21").Value = "  +10.39%  "
Set-TextCell "D22" "71.13"
$ws.Range("E22").Value = "  +5.32%  "
Set-TextCell "D23" "258.08"
$ws.Range("E23").Value = "  +4.99%  "
Set-TextCell "D24" "2.99"
$ws.Range("E24").Value = "  +5.17%  "
$ws.Range("E25").Value = "  +14.15%  "
Set-TextCell "D26" "28.05"
$ws.Range("E26").Value = "  +33.62%  "
$ws.Range("E27").Value = "  +0.04%  "
Set-TextCell "D28" "10.50"
$ws.Range("E28").Value = "  +7.30%  "
Set-TextCell "D29" "39.55"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  +3.64%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D31" "6.12"
$ws.Range("E31").Value = "  +10.86%  "
$ws.Range("B32").Value = "LidoDAOToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D32" "3.72"
$ws.Range("E32").Value = "  -1.43%  "
Set-TextCell "D33" "2.33"
$ws.Range("E33").Value = "  +22.64%  "
Set-TextCell "D34" "2.92"
$ws.Range("E34").Value = "  +5.01%  "
Set-TextCell "D35" "0.0836"
$ws.Range("E35").Value = "  +7.67%  "
Set-TextCell "D36" "150.11"
$ws.Range("E36").Value = "  +2.50%  "
$ws.Range("E37").Value = "  +4.96%  "
$ws.Range("E38").Value = "  +4.47%  "
$ws.Range("E39").Value = "  +6.61%  "
Set-TextCell "D40" "15.78"
$ws.Range("E40").Value = "  +5.44%  "
$ws.Range("E41").Value = "  +12.96%  "
Set-TextCell "D42" "0.0324"
$ws.Range("E42").Value = "  +7.63%  "
$ws.Range("D43").Value = "2.034.43"
$ws.Range("E43").Value = "  +6.45%  "
Set-TextCell "D44" "19.41"
$ws.Range("E44").Value = "  +34.19%  "
$ws.Range("E45").Value = "  -0.02%  "
Set-TextCell "D46" "91.64"
$ws.Range("E46").Value = "  -0.21%  "
Set-TextCell "D47" "1.79"
$ws.Range("E47").Value = "  -0.58%  "
Set-TextCell "D48" "9.16"
$ws.Range("E48").Value = "  +9.94%  "
Set-TextCell "D49" "109.10"
$ws.Range("E49").Value = "  +11.26%  "
$ws.Range("E50").Value = "  +7.59%  "
$ws.Range("D51").Value = "2.854.94"
$ws.Range("E51").Value = "  +10.14%  "
